$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-5 (A:AH). Row 6 is removed entirely.
$data = @(
    @(45136.50694444445, 5.237, 1.607, 0, 2.478, 3.053, 2.519, 5.488, 1.554, 0.9409999999999999, 4.022, 1.069, 0.9379999999999999, 0.6929999999999999, 0.87, 2.749, 1.106, 0.51, 0.063, 20.094, 4.803, 2.372, 3.891, 1.034, 0.249, 1.98, 1.154, 0.674, 0.9399999999999999, 3.129, 2.866, 3.128, 0.446, 1.43),
    @(45136.51388888889, 5.495, 3.344, 0, 8.061999999999999, 7.348, 3.882, 10.821, 4.539, 2.258, 4.397, 3.371, 3.292, 1.098, 2.877, 5.177, 2.541, 0.408, 0.034, 43.724, 9.314, 3.457, 6.484, 2.988, 0.444, 5.468, 2.668, 2.072, 2.541, 4.569, 1.133, 8.970000000000001, 1.485, 3.555),
    @(45136.52083333334, 15.15, 10.937, 0.349, 30.464, 25.666, 11.698, 39.661, 16.962, 7.796, 12.262, 12.356, 12.815, 3.627, 10.95, 16.421, 9.154999999999999, 0.41, 0.332, 164.077, 31.411, 10.634, 21.371, 10.968, 1.479, 20.071, 9.156000000000001, 7.835, 9.289999999999999, 13.646, 0.716, 35.501, 5.733, 12.787),
    @(45136.52777777778, 22.16, 16.34, 0.61, 46.37, 38.7, 17.3, 64.75, 25.79, 11.75, 17.96, 18.73, 19.58, 5.44, 16.68, 24.46, 13.87, 0.42, 0.55, 250.02, 47.29, 15.78, 32.1, 16.65, 2.22, 31.92, 13.77, 11.93, 14.09, 20.16, 0.52, 58.58, 8.75, 19.35)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $i + 2
    $rowVals = $data[$i]
    for ($c = 1; $c -le $rowVals.Length; $c++) {
        $ws.Cells.Item($rowIndex, $c).Value = $rowVals[$c - 1]
    }
}

# Delete row 6 entirely (was the 5th data row; sheet now ends at row 5)
$ws.Rows.Item(6).Delete()

# Narrow a handful of columns from width 8 to width 7.
# (ColumnWidth is in "characters"; the saved OOXML width = ColumnWidth + 5/6,
# so use 7 - 5/6 = 6.1666... to land on a stored width of exactly 7.)
$narrowWidth = 7 - (5 / 6)
$ws.Columns.Item(2).ColumnWidth = $narrowWidth
$ws.Columns.Item(15).ColumnWidth = $narrowWidth
$ws.Columns.Item(17).ColumnWidth = $narrowWidth
$ws.Columns.Item(27).ColumnWidth = $narrowWidth
$ws.Columns.Item(29).ColumnWidth = $narrowWidth
